$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Complete"
$ws.Range("D4").Value = "Complete"
$ws.Range("D22").Value = "Complete"
$ws.Range("C34").Value = "Complete"
$ws.Range("C35").Value = "Complete"
$ws.Range("C36").Value = "Complete"

$ws.PageSetup.PrintArea = '$A$1:$D$36'
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Zoom = 91

$ws.Range("B23").Select()
